$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from G1 (existing header cell) to H1, then set the new header text
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Values for new "Save" column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
